$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated distance values (columns C:J, rows 2:25) as captured in the target diff.
$newValues = @{
    "C2" = 0.1221089063380579
    "D2" = 0.05667183572656422
    "E2" = 0.3211595485804488
    "F2" = 0.2193388286095951
    "G2" = 0.1415814790529661
    "H2" = 0.1910585783618671
    "I2" = 0.02490447911794885
    "J2" = 0.08715245548174068
    "C3" = 0.0521406020907713
    "D3" = 0.3532589433180323
    "E3" = 0.6710434694221715
    "F3" = 0.5877335378547766
    "G3" = 0.5157913150014871
    "H3" = 0.5570805888616498
    "I3" = 0.01257038024501589
    "J3" = 0.3815832249538377
    "C4" = 0.06061850427318792
    "D4" = 0.05272820037210534
    "E4" = 0.2161304025387943
    "F4" = 0.3166335343525982
    "G4" = 0.2149017587941981
    "H4" = 0.2141689417392107
    "I4" = 0.1538739188402041
    "J4" = 0.06661064523310527
    "C5" = 0.05236527680284045
    "D5" = 0.2727019875170374
    "E5" = 0.4466589320811828
    "F5" = 0.5029279674534891
    "G5" = 0.5941951608508814
    "H5" = 0.5421197261852483
    "I5" = 0.07898644808681816
    "J5" = 0.2940878029595256
    "C6" = 0.2504702820923962
    "D6" = 0.1931191904018958
    "E6" = 0.2968538858314453
    "F6" = 0.2530993598845906
    "G6" = 0.1317844090981011
    "H6" = 0.1761312533493202
    "I6" = 0.04104805868853451
    "J6" = 0.2044926820048217
    "C7" = 0.1523271065783985
    "D7" = 0.4728680878765464
    "E7" = 0.3895165000223245
    "F7" = 0.3220876412301786
    "G7" = 0.3638931718703192
    "H7" = 0.3436376859969555
    "I7" = 0.03839212826994619
    "J7" = 0.4091561293709164
    "C8" = 0.1562042908344657
    "D8" = 0.1836026026006673
    "E8" = 0.4867276688979716
    "F8" = 0.3767857758787648
    "G8" = 0.3020545768112137
    "H8" = 0.3582755530718278
    "I8" = 0.196273976601093
    "J8" = 0.1801164855031319
    "C9" = 0.110999018595311
    "D9" = 0.05750540359501723
    "E9" = 0.2678921817048892
    "F9" = 0.2878747055431122
    "G9" = 0.2543634632311458
    "H9" = 0.300310610456655
    "I9" = 0.06591779100846847
    "J9" = 0.06450445058865004
    "C10" = 0.3555411626094466
    "D10" = 0.3890580211478373
    "E10" = 0.1976322408940027
    "F10" = 0.06426727500536926
    "G10" = 0.1315463511749019
    "H10" = 0.07470209583380527
    "I10" = 0.1249252173210046
    "J10" = 0.2191580708406993
    "C11" = 0.2916408856922709
    "D11" = 0.2949881489384521
    "E11" = 0.3486749359084462
    "F11" = 0.3783540089594553
    "G11" = 0.1657765384589845
    "H11" = 0.3133925709457292
    "I11" = 0.4640262117893734
    "J11" = 0.2180542334657396
    "C12" = 0.05694529914590994
    "D12" = 0.1261039145409353
    "E12" = 0.3993798401154126
    "F12" = 0.2232839199634121
    "G12" = 0.209771203020897
    "H12" = 0.2251856991234587
    "I12" = 0.1045951669257556
    "J12" = 0.1073151451188274
    "C13" = 0.02999454772084903
    "D13" = 0.1655854160283751
    "E13" = 0.2850121201910394
    "F13" = 0.2024456846945018
    "G13" = 0.1743710317806629
    "H13" = 0.1475431396275167
    "I13" = 0.2187952918030596
    "J13" = 0.1387895025337578
    "C14" = 0.3346598016955187
    "D14" = 0.4173411504166115
    "E14" = 0.5155403708011996
    "F14" = 0.5767381663888015
    "G14" = 0.4962656372070536
    "H14" = 0.4904983353506448
    "I14" = 0.02191923626735155
    "J14" = 0.3340372793447789
    "C15" = 0.3101110215720707
    "D15" = 0.2919151534546408
    "E15" = 0.4034107821138355
    "F15" = 0.300149119480917
    "G15" = 0.2029541637292397
    "H15" = 0.3171388721944635
    "I15" = 0.03392691159027722
    "J15" = 0.302848360067093
    "C16" = 0.1430306054787642
    "D16" = 0.2708874560637792
    "E16" = 0.6177783880851135
    "F16" = 0.4660143435166007
    "G16" = 0.3049137444030126
    "H16" = 0.333223401824376
    "I16" = 0.1827889953002473
    "J16" = 0.1243117779858766
    "C17" = 0.2986124684659847
    "D17" = 0.2915008351073351
    "E17" = 0.6204184786470422
    "F17" = 0.6651269427960155
    "G17" = 0.533512656677863
    "H17" = 0.5879749931444529
    "I17" = 0.006189272009205444
    "J17" = 0.2911961486013736
    "C18" = 0.3049196220952002
    "D18" = 0.5867896226282683
    "E18" = 0.6646755233824322
    "F18" = 0.6527563095190463
    "G18" = 0.5702477869729837
    "H18" = 0.488614381560924
    "I18" = 0.06572905670118356
    "J18" = 0.4061570297407783
    "C19" = 0.2730433218116961
    "D19" = 0.3692725404981654
    "E19" = 0.3643682281803224
    "F19" = 0.5218528830593737
    "G19" = 0.4739422182336951
    "H19" = 0.5146898778335515
    "I19" = 0.01715176490101544
    "J19" = 0.3073278356572086
    "C20" = 0.2025036683650922
    "D20" = 0.214333569215466
    "E20" = 0.2973776141360639
    "F20" = 0.2403703570472182
    "G20" = 0.2619818583713167
    "H20" = 0.2251790251343379
    "I20" = 0.1117418270225956
    "J20" = 0.1780892119188407
    "C21" = 0.2295611529714554
    "D21" = 0.4221420735844187
    "E21" = 0.6061319196989261
    "F21" = 0.5387636084583555
    "G21" = 0.5472891506367564
    "H21" = 0.5075474056752981
    "I21" = 0.1732494065400206
    "J21" = 0.3944532975639055
    "C22" = 0.2720063863017567
    "D22" = 0.1567008811481451
    "E22" = 0.427342897950528
    "F22" = 0.2864239863294233
    "G22" = 0.1418987071935987
    "H22" = 0.2245200748804823
    "I22" = 0.0446402282581856
    "J22" = 0.1133086345571219
    "C23" = 0.150641714080226
    "D23" = 0.3833547985668339
    "E23" = 0.6234743306603664
    "F23" = 0.6558536924011198
    "G23" = 0.5177023903577762
    "H23" = 0.5753083671493943
    "I23" = 0.04562550088336464
    "J23" = 0.3367232312751348
    "C24" = 0.1540082716866467
    "D24" = 0.3427044964394886
    "E24" = 0.5456800975731519
    "F24" = 0.515025256662373
    "G24" = 0.4655227203131477
    "H24" = 0.4914646086840925
    "I24" = 0.1335068387041827
    "J24" = 0.2863582778424957
    "C25" = 0.3323439558694475
    "D25" = 0.4323057290783849
    "E25" = 0.2861082711539883
    "F25" = 0.475748258487646
    "G25" = 0.4547685833270968
    "H25" = 0.391163424321444
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

Write-Output "Updated $($newValues.Count) cells"
